# Append a new "publicform_admin_error" mail-template row (row 14) to the
# mail_template sheet, mirroring the existing rows' column layout:
#   A=id  B=mail_key_name  C=mail_view_name  D=mail_template_type
#   E=mail_subject  F=mail_body
#
# Shared strings must be created in B, E, C, F, D order so the appended
# <si> entries in sharedStrings.xml land in the same sequence as the
# target workbook (D14 reuses the existing "body" string, so its position
# doesn't add a new shared string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = 'publicform_admin_error'
$ws.Range("E14").Value = '[${system:site_name}]公開フォーム ${publicform:public_form_view_name} でエラーが発生しました'
$ws.Range("C14").Value = 'Public form administrator error'
$ws.Range("F14").Value = '${user:user_name}\n\nAn error occurred while filling out the public form ${publicform:public_form_view_name}.\nPlease check the entered contents and the error contents.\n\n■ Input contents\n${form:values}\n\n■ Error summary\n${error:message}\n\n■ Error details\n${error:stacktrace}'
$ws.Range("D14").Value = 'body'

# Matches the author's final selection landing on the new subject cell.
$ws.Range("F14").Select() | Out-Null
